$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3749.875
$ws.Range("I64").Value = 3599.8
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3599.8
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3351.8
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 3749.875
$ws.Range("I67").Value = 3599.8
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3599.8
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2741.8
$ws.Range("N67").Value = -5716

$ws.Range("H69").Value = 15658
$ws.Range("I69").Value = 10966
$ws.Range("J69").Value = 17417.5
$ws.Range("K69").Value = 32898
$ws.Range("L69").Value = 52252.5
$ws.Range("M69").Value = -32024
$ws.Range("N69").Value = -54000.5

$ws.Range("H72").Value = 15658
$ws.Range("I72").Value = 10966
$ws.Range("J72").Value = 17417.5
$ws.Range("K72").Value = 98694
$ws.Range("L72").Value = 156757.5
$ws.Range("M72").Value = -94326
$ws.Range("N72").Value = -165493.5

$ws.Range("H97").Value = 3425.4285
$ws.Range("I97").Value = 1990
$ws.Range("J97").Value = 3497.2
$ws.Range("K97").Value = 5970
$ws.Range("L97").Value = 10491.6
$ws.Range("M97").Value = -5474
$ws.Range("N97").Value = -11483.6

$ws.Range("H112").Value = 3349.923
$ws.Range("J112").Value = 3187.0571
$ws.Range("L112").Value = 9561.1713
$ws.Range("N112").Value = -11777.1713

$ws.Range("H132").Value = 5947.5415
$ws.Range("I132").Value = 4949
$ws.Range("J132").Value = 10940.25
$ws.Range("K132").Value = 14847
$ws.Range("L132").Value = 32820.75
$ws.Range("M132").Value = -12317
$ws.Range("N132").Value = -37880.75

$ws.Range("H137").Value = 6683.857
$ws.Range("I137").Value = 7598.227
$ws.Range("K137").Value = 22794.681
$ws.Range("M137").Value = -20244.681

$ws.Range("H138").Value = 8476.75
$ws.Range("I138").Value = 7613.1
$ws.Range("J138").Value = 8572.710999999999
$ws.Range("K138").Value = 22839.3
$ws.Range("L138").Value = 25718.133
$ws.Range("M138").Value = -17699.3
$ws.Range("N138").Value = -35998.133

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15324.268
$ws.Range("I32").Value = 14829.028
$ws.Range("K32").Value = 14829.028
$ws.Range("M32").Value = -14542.028

$ws.Range("H74").Value = 1249.7949
$ws.Range("I74").Value = 982.2432
$ws.Range("K74").Value = 982.2432
$ws.Range("M74").Value = -108.2432

$ws.Range("H77").Value = 1249.7949
$ws.Range("I77").Value = 982.2432
$ws.Range("K77").Value = 4911.216
$ws.Range("M77").Value = -543.2160000000003

$ws.Range("H133").Value = 129999.5
$ws.Range("J133").Value = 129999.5
$ws.Range("L133").Value = 129999.5
$ws.Range("N133").Value = -135059.5

$ws.Range("H138").Value = 124718.8
$ws.Range("J138").Value = 124718.8
$ws.Range("L138").Value = 124718.8
$ws.Range("N138").Value = -134998.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3314.1
$ws.Range("I86").Value = 2625.6667
$ws.Range("K86").Value = 2625.6667
$ws.Range("M86").Value = -1502.6667

$ws.Range("H89").Value = 3314.1
$ws.Range("I89").Value = 2625.6667
$ws.Range("K89").Value = 13128.3335
$ws.Range("M89").Value = -7512.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4726.3794
$ws.Range("I31").Value = 1473.0555
$ws.Range("J31").Value = 10050
$ws.Range("K31").Value = 1473.0555
$ws.Range("L31").Value = 10050
$ws.Range("M31").Value = -1178.0555
$ws.Range("N31").Value = -10640

$ws.Range("H34").Value = 4726.3794
$ws.Range("I34").Value = 1473.0555
$ws.Range("J34").Value = 10050
$ws.Range("K34").Value = 1473.0555
$ws.Range("L34").Value = 10050
$ws.Range("M34").Value = -1271.0555
$ws.Range("N34").Value = -10454

$ws.Range("H58").Value = 2951.4614
$ws.Range("I58").Value = 3176.8
$ws.Range("K58").Value = 3176.8
$ws.Range("M58").Value = -2973.8

$ws.Range("H62").Value = 166673470
$ws.Range("J62").Value = 6833
$ws.Range("L62").Value = 6833
$ws.Range("N62").Value = -8081

$ws.Range("H65").Value = 166673470
$ws.Range("J65").Value = 6833
$ws.Range("L65").Value = 34165
$ws.Range("N65").Value = -40405

$ws.Range("H132").Value = 6185.2334
$ws.Range("J132").Value = 4752.8184
$ws.Range("L132").Value = 14258.4552
$ws.Range("N132").Value = -19318.4552

$ws.Range("H134").Value = 4210.1304
$ws.Range("I134").Value = 4223.7295
$ws.Range("J134").Value = 4154.222
$ws.Range("K134").Value = 12671.1885
$ws.Range("L134").Value = 12462.666
$ws.Range("M134").Value = -10136.1885
$ws.Range("N134").Value = -17532.666

$ws.Range("H136").Value = 2951.4614
$ws.Range("I136").Value = 3176.8
$ws.Range("K136").Value = 9530.400000000001
$ws.Range("M136").Value = -6980.400000000001

$ws.Range("H141").Value = 321729.56
$ws.Range("I141").Value = 117499.5
$ws.Range("J141").Value = 347258.3
$ws.Range("K141").Value = 117499.5
$ws.Range("L141").Value = 347258.3
$ws.Range("M141").Value = -112319.5
$ws.Range("N141").Value = -357618.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1578.65
$ws.Range("I5").Value = 1165.3103
$ws.Range("K5").Value = 3495.9309
$ws.Range("M5").Value = -3383.9309

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H34").Value = 665.875
$ws.Range("I34").Value = 587.8333
$ws.Range("J34").Value = 900
$ws.Range("K34").Value = 1763.4999
$ws.Range("L34").Value = 2700
$ws.Range("M34").Value = -1679.4999
$ws.Range("N34").Value = -2868

$ws.Range("H39").Value = 1562
$ws.Range("J39").Value = 2193
$ws.Range("L39").Value = 6579
$ws.Range("N39").Value = -7167

$ws.Range("H55").Value = 14775.857
$ws.Range("J55").Value = 99999
$ws.Range("L55").Value = 299997
$ws.Range("N55").Value = -300351

$ws.Range("H121").Value = 14739.027
$ws.Range("J121").Value = 22902.521
$ws.Range("L121").Value = 68707.56299999999
$ws.Range("N121").Value = -71327.56299999999

$ws.Range("H128").Value = 183961.58
$ws.Range("I128").Value = 183961.58
$ws.Range("K128").Value = 551884.74
$ws.Range("M128").Value = -546904.74

$ws.Range("H129").Value = 1025.0769
$ws.Range("J129").Value = 1847.5
$ws.Range("L129").Value = 5542.5
$ws.Range("N129").Value = -15542.5

$ws.Range("H135").Value = 1578.65
$ws.Range("I135").Value = 1165.3103
$ws.Range("K135").Value = 10487.7927
$ws.Range("M135").Value = -7952.792700000002

$ws.Range("H137").Value = 8159.3335
$ws.Range("J137").Value = 6000
$ws.Range("L137").Value = 18000
$ws.Range("N137").Value = -28200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2710.524
$ws.Range("I107").Value = 2594.0908
$ws.Range("K107").Value = 2594.0908
$ws.Range("M107").Value = -674.0907999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9998.4
$ws.Range("I7").Value = 9998
$ws.Range("K7").Value = 9998
$ws.Range("M7").Value = -9886

$ws.Range("H20").Value = 13000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H22").Value = 1211.8667
$ws.Range("I22").Value = 1336.75
$ws.Range("K22").Value = 1336.75
$ws.Range("M22").Value = -1041.75

$ws.Range("H27").Value = 1211.8667
$ws.Range("I27").Value = 1336.75
$ws.Range("K27").Value = 1336.75
$ws.Range("M27").Value = -1229.75

$ws.Range("H68").Value = 2469.25
$ws.Range("I68").Value = 2343.1
$ws.Range("J68").Value = 3100
$ws.Range("K68").Value = 2343.1
$ws.Range("L68").Value = 3100
$ws.Range("M68").Value = -1594.1
$ws.Range("N68").Value = -4598

$ws.Range("H71").Value = 2469.25
$ws.Range("I71").Value = 2343.1
$ws.Range("J71").Value = 3100
$ws.Range("K71").Value = 11715.5
$ws.Range("L71").Value = 15500
$ws.Range("M71").Value = -7971.5
$ws.Range("N71").Value = -22988

$ws.Range("H126").Value = 9998.4
$ws.Range("I126").Value = 9998
$ws.Range("K126").Value = 29994
$ws.Range("M126").Value = -27524

$ws.Range("H139").Value = 95638.64
$ws.Range("J139").Value = 93540.09
$ws.Range("L139").Value = 93540.09
$ws.Range("N139").Value = -103820.09

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4437.3184
$ws.Range("I136").Value = 3971.8276
$ws.Range("K136").Value = 11915.4828
$ws.Range("M136").Value = -9365.4828
